$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.62515252014563
$ws.Range("C2").Value = 9.811478452150341
$ws.Range("D2").Value = 9.917674862594623
$ws.Range("F2").Value = 29.46985944617
$ws.Range("G2").Value = 28.72111806352111
$ws.Range("H2").Value = 14.28172866198689
$ws.Range("I2").Value = 22.0327378918661
$ws.Range("J2").Value = 10.19790059459565
$ws.Range("L2").Value = 11.78708295704407
$ws.Range("N2").Value = 17.31632836771244
$ws.Range("O2").Value = 21.73491782278565

$ws.Range("B3").Value = 15.14697405513095
$ws.Range("C3").Value = 9.649048914558024
$ws.Range("D3").Value = 9.919213554014245
$ws.Range("F3").Value = 29.52078107744305
$ws.Range("G3").Value = 28.73729906796211
$ws.Range("H3").Value = 14.32272593277829
$ws.Range("I3").Value = 22.12957570788443
$ws.Range("J3").Value = 10.22257527959585
$ws.Range("L3").Value = 11.76818305355794
$ws.Range("N3").Value = 17.3484877459503
$ws.Range("O3").Value = 21.79220438846713

$ws.Range("B4").Value = 14.84678542920077
$ws.Range("C4").Value = 9.54730978386425
$ws.Range("D4").Value = 9.921465334789563
$ws.Range("F4").Value = 29.55895262571818
$ws.Range("G4").Value = 28.75650776443739
$ws.Range("H4").Value = 14.35019841037687
$ws.Range("I4").Value = 22.193230952739
$ws.Range("J4").Value = 10.23857906140492
$ws.Range("L4").Value = 11.75805213330639
$ws.Range("N4").Value = 17.36991327984076
$ws.Range("O4").Value = 21.8320993693921

$ws.Range("B5").Value = 14.72297732633399
$ws.Range("C5").Value = 9.505380163925988
$ws.Range("D5").Value = 9.922712608509702
$ws.Range("F5").Value = 29.57624174232339
$ws.Range("G5").Value = 28.76666346695142
$ws.Range("H5").Value = 14.36197187118806
$ws.Range("I5").Value = 22.22022579768477
$ws.Range("J5").Value = 10.24531585858988
$ws.Range("L5").Value = 11.75429739411953
$ws.Range("N5").Value = 17.37906738943455
$ws.Range("O5").Value = 21.8495418409848

$ws.Range("B6").Value = 14.70233521026688
$ws.Range("C6").Value = 9.498390420132445
$ws.Range("D6").Value = 9.92293965497819
$ws.Range("F6").Value = 29.57921722583015
$ws.Range("G6").Value = 28.76849027288688
$ws.Range("H6").Value = 14.36396176284162
$ws.Range("I6").Value = 22.22477196799108
$ws.Range("J6").Value = 10.24644750720768
$ws.Range("L6").Value = 11.75369657972793
$ws.Range("N6").Value = 17.38061299285594
$ws.Range("O6").Value = 21.85250965170732

$ws.Range("B7").Value = 14.84512145348212
$ws.Range("C7").Value = 9.54674616333979
$ws.Range("D7").Value = 9.92148081994535
$ws.Range("F7").Value = 29.55917877603673
$ws.Range("G7").Value = 28.75663530818447
$ws.Range("H7").Value = 14.35035485047229
$ws.Range("I7").Value = 22.19359074427107
$ws.Range("J7").Value = 10.23866904446871
$ws.Range("L7").Value = 11.75799997854266
$ws.Range("N7").Value = 17.37003502160492
$ws.Range("O7").Value = 21.83232980983914

$ws.Range("B8").Value = 15.46174533901936
$ws.Range("C8").Value = 9.755909590654788
$ws.Range("D8").Value = 9.917934753424465
$ws.Range("F8").Value = 29.48598262740016
$ws.Range("G8").Value = 28.72477071747652
$ws.Range("H8").Value = 14.2953870753349
$ws.Range("I8").Value = 22.06525635063181
$ws.Range("J8").Value = 10.20623161421862
$ws.Range("L8").Value = 11.78026233546459
$ws.Range("N8").Value = 17.32706874801295
$ws.Range("O8").Value = 21.75368919420436

$ws.Range("B9").Value = 16.61161069037694
$ws.Range("C9").Value = 10.1487023524569
$ws.Range("D9").Value = 9.9213050711687
$ws.Range("F9").Value = 29.39732950394592
$ws.Range("G9").Value = 28.73596443484823
$ws.Range("H9").Value = 14.20585101675016
$ws.Range("I9").Value = 21.84690817395827
$ws.Range("J9").Value = 10.14936947550584
$ws.Range("L9").Value = 11.83546726764621
$ws.Range("N9").Value = 17.25610807117304
$ws.Range("O9").Value = 21.63702183898411

$ws.Range("B10").Value = 17.41173732587275
$ws.Range("C10").Value = 10.42484340890823
$ws.Range("D10").Value = 9.930008642054872
$ws.Range("F10").Value = 29.36575046838841
$ws.Range("G10").Value = 28.78913071658335
$ws.Range("H10").Value = 14.15120583848246
$ws.Range("I10").Value = 21.70682270755463
$ws.Range("J10").Value = 10.1116727739461
$ws.Range("L10").Value = 11.88286188722982
$ws.Range("N10").Value = 17.21203824252591
$ws.Range("O10").Value = 21.57430796891384

$ws.Range("B11").Value = 17.7645062213983
$ws.Range("C11").Value = 10.54740178290301
$ws.Range("D11").Value = 9.935304257001807
$ws.Range("F11").Value = 29.35867857035989
$ws.Range("G11").Value = 28.8230436434345
$ws.Range("H11").Value = 14.1287660011363
$ws.Range("I11").Value = 21.64751586184641
$ws.Range("J11").Value = 10.09540237206673
$ws.Range("L11").Value = 11.90585943253437
$ws.Range("N11").Value = 17.19373246646055
$ws.Range("O11").Value = 21.55079226641928

$ws.Range("B12").Value = 17.89636818299537
$ws.Range("C12").Value = 10.59334310849499
$ws.Range("D12").Value = 9.937500180918784
$ws.Range("F12").Value = 29.35704896620915
$ws.Range("G12").Value = 28.83727851281085
$ws.Range("H12").Value = 14.12061650434213
$ws.Range("I12").Value = 21.62569399760605
$ws.Range("J12").Value = 10.08936692181964
$ws.Range("L12").Value = 11.91477031577785
$ws.Range("N12").Value = 17.18705034541754
$ws.Range("O12").Value = 21.54260961295095

$ws.Range("B13").Value = 17.86804767296143
$ws.Range("C13").Value = 10.583470118336
$ws.Range("D13").Value = 9.937018798995242
$ws.Range("F13").Value = 29.35735331811041
$ws.Range("G13").Value = 28.8341509503686
$ws.Range("H13").Value = 14.12235616472952
$ws.Range("I13").Value = 21.63036540926288
$ws.Range("J13").Value = 10.09066117639239
$ws.Range("L13").Value = 11.91284227786922
$ws.Range("N13").Value = 17.18847835567842
$ws.Range("O13").Value = 21.54433975343474

$ws.Range("B14").Value = 17.77538973148496
$ws.Range("C14").Value = 10.55119096013305
$ws.Range("D14").Value = 9.935481106647197
$ws.Range("F14").Value = 29.35852349581011
$ws.Range("G14").Value = 28.82418688240139
$ws.Range("H14").Value = 14.12808856197591
$ws.Range("I14").Value = 21.64570780499986
$ws.Range("J14").Value = 10.09490331338145
$ws.Range("L14").Value = 11.90658851008038
$ws.Range("N14").Value = 17.19317771976465
$ws.Range("O14").Value = 21.55010459312424

$ws.Range("B15").Value = 17.71840645029598
$ws.Range("C15").Value = 10.53135715240857
$ws.Range("D15").Value = 9.93456399769544
$ws.Range("F15").Value = 29.35937676788467
$ws.Range("G15").Value = 28.81826477381746
$ws.Range("H15").Value = 14.1316451439669
$ws.Range("I15").Value = 21.65518836563576
$ws.Range("J15").Value = 10.09751811381116
$ws.Range("L15").Value = 11.90278409391737
$ws.Range("N15").Value = 17.19608874010895
$ws.Range("O15").Value = 21.55372981573859

$ws.Range("B16").Value = 17.38844760373885
$ws.Range("C16").Value = 10.41676985670136
$ws.Range("D16").Value = 9.929689319279625
$ws.Range("F16").Value = 29.36635937666583
$ws.Range("G16").Value = 28.78710969628294
$ws.Range("H16").Value = 14.15272101935209
$ws.Range("I16").Value = 21.7107875454286
$ws.Range("J16").Value = 10.11275371147615
$ws.Range("L16").Value = 11.88138748955499
$ws.Range("N16").Value = 17.21326956710024
$ws.Range("O16").Value = 21.57594577354599

$ws.Range("B17").Value = 17.18307341033708
$ws.Range("C17").Value = 10.34566955514293
$ws.Range("D17").Value = 9.927040039076097
$ws.Range("F17").Value = 29.3725109384863
$ws.Range("G17").Value = 28.77048480088761
$ws.Range("H17").Value = 14.16626998019048
$ws.Range("I17").Value = 21.74602824520434
$ws.Range("J17").Value = 10.12232480934733
$ws.Range("L17").Value = 11.86862640430399
$ws.Range("N17").Value = 17.22425514553887
$ws.Range("O17").Value = 21.59085956811045

$ws.Range("B18").Value = 17.06389808315565
$ws.Range("C18").Value = 10.3044881880274
$ws.Range("D18").Value = 9.925642124122531
$ws.Range("F18").Value = 29.37673572675277
$ws.Range("G18").Value = 28.76183866157474
$ws.Range("H18").Value = 14.17429061689284
$ws.Range("I18").Value = 21.7667136174123
$ws.Range("J18").Value = 10.12791252194152
$ws.Range("L18").Value = 11.86142215112896
$ws.Range("N18").Value = 17.2307377421856
$ws.Range("O18").Value = 21.59990938435843

$ws.Range("B19").Value = 17.02337070851666
$ws.Range("C19").Value = 10.29049661016959
$ws.Range("D19").Value = 9.925190480519086
$ws.Range("F19").Value = 29.37828408937566
$ws.Range("G19").Value = 28.75906871334095
$ws.Range("H19").Value = 14.17704535629829
$ws.Range("I19").Value = 21.77378872416345
$ws.Range("J19").Value = 10.1298186382129
$ws.Range("L19").Value = 11.85900633560512
$ws.Range("N19").Value = 17.23296082025021
$ws.Range("O19").Value = 21.60305448305452

$ws.Range("B20").Value = 17.20504531228905
$ws.Range("C20").Value = 10.35326815469288
$ws.Range("D20").Value = 9.927309042336224
$ws.Range("F20").Value = 29.37178503757919
$ws.Range("G20").Value = 28.77215977793783
$ws.Range("H20").Value = 14.16480410950129
$ws.Range("I20").Value = 21.74223376652397
$ws.Range("J20").Value = 10.12129739692917
$ws.Range("L20").Value = 11.86997084299042
$ws.Range("N20").Value = 17.22306874467977
$ws.Range("O20").Value = 21.58922312955711

$ws.Range("B21").Value = 17.80265323241105
$ws.Range("C21").Value = 10.56068507876803
$ws.Range("D21").Value = 9.935927604646068
$ws.Range("F21").Value = 29.35815134056537
$ws.Range("G21").Value = 28.82707582822937
$ws.Range("H21").Value = 14.12639537231999
$ws.Range("I21").Value = 21.6411840934024
$ws.Range("J21").Value = 10.09365388443825
$ws.Range("L21").Value = 11.90841994209979
$ws.Range("N21").Value = 17.19179062566509
$ws.Range("O21").Value = 21.54839170924907

$ws.Range("B22").Value = 18.18313255433409
$ws.Range("C22").Value = 10.69349882137723
$ws.Range("D22").Value = 9.94267043429214
$ws.Range("F22").Value = 29.35535096009215
$ws.Range("G22").Value = 28.87108071290964
$ws.Range("H22").Value = 14.10332142684289
$ws.Range("I22").Value = 21.57885189295406
$ws.Range("J22").Value = 10.07632028748427
$ws.Range("L22").Value = 11.93472507664307
$ws.Range("N22").Value = 17.17280478141075
$ws.Range("O22").Value = 21.52591619157851

$ws.Range("B23").Value = 17.98102078683612
$ws.Range("C23").Value = 10.62287394742336
$ws.Range("D23").Value = 9.938970630145336
$ws.Range("F23").Value = 29.35628682445017
$ws.Range("G23").Value = 28.84685450970998
$ws.Range("H23").Value = 14.11545076842437
$ws.Range("I23").Value = 21.61178001244976
$ws.Range("J23").Value = 10.08550463204994
$ws.Range("L23").Value = 11.92057942350623
$ws.Range("N23").Value = 17.18280483443605
$ws.Range("O23").Value = 21.53752618529346

$ws.Range("B24").Value = 17.1951152423013
$ws.Range("C24").Value = 10.34983377520438
$ws.Range("D24").Value = 9.927187035897326
$ws.Range("F24").Value = 29.37211107375382
$ws.Range("G24").Value = 28.77139968040792
$ws.Range("H24").Value = 14.16546610960848
$ws.Range("I24").Value = 21.74394792648143
$ws.Range("J24").Value = 10.12176162478115
$ws.Range("L24").Value = 11.86936260993619
$ws.Range("N24").Value = 17.22360459676354
$ws.Range("O24").Value = 21.58996148183053

$ws.Range("B25").Value = 16.30782542469303
$ws.Range("C25").Value = 10.04449831112532
$ws.Range("D25").Value = 9.919294510368621
$ws.Range("F25").Value = 29.41542453951644
$ws.Range("G25").Value = 28.72503851612235
$ws.Range("H25").Value = 14.22811761932554
$ws.Range("I25").Value = 21.90240806369775
$ws.Range("J25").Value = 10.16403330919563
$ws.Range("L25").Value = 11.8193168204338
$ws.Range("N25").Value = 17.27388560070303
$ws.Range("O25").Value = 21.66455141250041
